$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of test data (Tester list table grows by one row)
$ws.Range("A7").Value = "admin@devtest.com"
$ws.Range("B7").Value = "admin@devtest.com"
$ws.Range("C7").Value = "secret"

# Update view: zoom to 200% and select the whole table range
$ws.Application.ActiveWindow.Zoom = 200
$ws.Range("A1:C7").Select()
